$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
}

$ws.Range("D2").Value = "71.470.78"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.825.85"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue "D5" "703.35"
$ws.Range("E5").Value = "  -1.66%  "
Set-TextValue "D6" "171.48"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "3.826.49"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  +0.12%  "
Set-TextValue "D9" "0.527"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  -1.93%  "
Set-TextValue "D11" "7.49"
$ws.Range("E11").Value = "  +1.72%  "
Set-TextValue "D12" "0.489"
$ws.Range("E12").Value = "  +5.89%  "
$ws.Range("E13").Value = "  -1.71%  "
Set-TextValue "D14" "36.83"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "4.466.50"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "3.828.71"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "71.492.26"
$ws.Range("E17").Value = "  +0.22%  "
Set-TextValue "D18" "7.26"
$ws.Range("E18").Value = "  -0.16%  "
Set-TextValue "D19" "17.67"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("E20").Value = "  +0.15%  "
Set-TextValue "D21" "514.23"
$ws.Range("E21").Value = "  +2.94%  "
Set-TextValue "D22" "10.56"
$ws.Range("E22").Value = "  -1.56%  "
Set-TextValue "D23" "0.719"
$ws.Range("E23").Value = "  -2.60%  "
Set-TextValue "D24" "83.93"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  -2.73%  "
Set-TextValue "D26" "12.80"
$ws.Range("E26").Value = "  +4.77%  "
$ws.Range("D27").Value = "3.965.39"
$ws.Range("E27").Value = "  -0.72%  "
Set-TextValue "D28" "10.39"
$ws.Range("E28").Value = "  -2.86%  "
$ws.Range("E29").Value = "  +0.07%  "
Set-TextValue "D30" "2.01"
$ws.Range("E30").Value = "  -4.17%  "
$ws.Range("E31").Value = "  -5.13%  "
Set-TextValue "D32" "2.28"
$ws.Range("E32").Value = "  +1.34%  "
Set-TextValue "D33" "7.41"
$ws.Range("E33").Value = "  -1.27%  "
Set-TextValue "D34" "29.38"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D37").Value = "3.788.92"
$ws.Range("E37").Value = "  -0.50%  "
Set-TextValue "D38" "0.999"
$ws.Range("E38").Value = "  +0.16%  "
Set-TextValue "D39" "6.68"
$ws.Range("E39").Value = "  +10.75%  "
Set-TextValue "D40" "0.102"
$ws.Range("E40").Value = "  -1.60%  "
Set-TextValue "D41" "2.41"
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("E42").Value = "  -2.06%  "
Set-TextValue "D43" "3.25"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("E45").Value = "  -0.06%  "
Set-TextValue "D46" "166.36"
$ws.Range("E46").Value = "  +1.67%  "
Set-TextValue "D47" "49.99"
$ws.Range("E47").Value = "  +2.11%  "
Set-TextValue "D48" "432.98"
$ws.Range("E48").Value = "  +1.34%  "
Set-TextValue "D49" "0.000304"
$ws.Range("E49").Value = "  -5.73%  "
Set-TextValue "D50" "30.94"
$ws.Range("E50").Value = "  +9.21%  "
Set-TextValue "D51" "8.69"
$ws.Range("E51").Value = "  -0.62%  "

# Rows 35/36: Kaspa and Aptos swap positions with updated values
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D35" "9.37"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D36" "0.172"
$ws.Range("E36").Value = "  -5.68%  "
